$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 865: repeated column-header row, matching the periodic header
# rows already present every 36 rows in this time-series sheet (e.g.
# row 829). Values first, then clone the formatting (bold + border +
# centered alignment) from the previous header block.
# ------------------------------------------------------------------
$headerRow = 865
$headers = @("States/UT", "Active Cases", "Active Cases Since Yesterday", "Recovered Cases", "Recovered Cases Since Yesterday", "Deceased Cases", "Deceased Cases Since Yesterday", "Date")
for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item($headerRow, $col).Value = $headers[$col - 1]
}
$ws.Range("A829:H829").Copy()
$ws.Range("A865:H865").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Rows 866-900: the 11-10-2020 snapshot, one row per state/UT, same
# column layout as the rest of the sheet.
# ------------------------------------------------------------------
$data = @(
    @("Andaman and Nicobar Islands", 193, -4, 3744, 20, 55, 0, "11-10-2020"),
    @("Andhra Pradesh", 46624, -1041, 697699, 6659, 6194, 35, "11-10-2020"),
    @("Arunachal Pradesh", 2940, 80, 9035, 158, 23, 0, "11-10-2020"),
    @("Assam", 29221, -489, 163355, 1451, 811, 9, "11-10-2020"),
    @("Bihar", 11165, -109, 183390, 1269, 944, 10, "11-10-2020"),
    @("Chandigarh", 1229, -63, 11662, 157, 190, 2, "11-10-2020"),
    @("Chhattisgarh", 27369, -70, 111654, 2719, 1235, 39, "11-10-2020"),
    @("Dadra and Nagar Haveli and Daman and Diu", 100, -1, 3050, 13, 2, 0, "11-10-2020"),
    @("Delhi", 22007, 52, 278812, 2766, 5740, 48, "11-10-2020"),
    @("Goa", 4658, -125, 32777, 460, 499, 8, "11-10-2020"),
    @("Gujarat", 15936, -245, 130760, 1456, 3557, 10, "11-10-2020"),
    @("Haryana", 10677, -153, 128841, 1301, 1572, 10, "11-10-2020"),
    @("Himachal Pradesh", 2718, -138, 14278, 402, 248, 3, "11-10-2020"),
    @("Jammu and Kashmir", 10796, -348, 70955, 976, 1313, 7, "11-10-2020"),
    @("Jharkhand", 8362, -457, 82805, 1151, 784, 3, "11-10-2020"),
    @("Karnataka", 120948, 2078, 569947, 8337, 9891, 102, "11-10-2020"),
    @("Kerala", 96003, 4162, 182874, 7570, 978, 23, "11-10-2020"),
    @("Ladakh", 1022, -16, 3973, 87, 64, 1, "11-10-2020"),
    @("Madhya Pradesh", 15612, -556, 127034, 2147, 2599, 25, "11-10-2020"),
    @("Maharashtra", 221615, -15332, 1255779, 26440, 40040, 308, "11-10-2020"),
    @("Manipur", 2608, -250, 10396, 530, 88, 2, "11-10-2020"),
    @("Meghalaya", 2437, 13, 5045, 142, 62, 1, "11-10-2020"),
    @("Mizoram", 191, 3, 1984, 10, 0, 0, "11-10-2020"),
    @("Nagaland", 1238, 23, 5694, 38, 17, 0, "11-10-2020"),
    @("Odisha", 24414, -1046, 224273, 3885, 1006, 15, "11-10-2020"),
    @("Puducherry", 4719, -84, 25955, 412, 559, 1, "11-10-2020"),
    @("Punjab", 9752, -401, 109767, 1234, 3798, 25, "11-10-2020"),
    @("Rajasthan", 21354, -44, 133918, 2152, 1636, 15, "11-10-2020"),
    @("Sikkim", 450, -54, 2816, 95, 55, 2, "11-10-2020"),
    @("Tamil Nadu", 44150, -47, 597033, 5222, 10187, 67, "11-10-2020"),
    @("Telengana", 25713, -391, 185128, 2103, 1222, 5, "11-10-2020"),
    @("Tripura", 3951, -88, 24086, 285, 315, 2, "11-10-2020"),
    @("Uttarakhand", 7321, 32, 46470, 412, 734, 18, "11-10-2020"),
    @("Uttar Pradesh", 40210, -1077, 387149, 4063, 6353, 60, "11-10-2020"),
    @("West Bengal", 29793, 497, 255838, 3032, 5563, 62, "11-10-2020")
)

$startRow = 866
$endRow = $startRow + $data.Count - 1

# Column H holds dates written as plain "DD-MM-YYYY" text (like the rest
# of the sheet, e.g. H2 = "13-09-2020"), not real date serials, so format
# the column as Text first to stop the day-month-year string being
# auto-parsed into a date value.
$ws.Range("H" + $startRow + ":H" + $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
